$d = $word.ActiveDocument

# --- Merge split text runs in Title / Author / Abstract paragraphs ---
$d.Paragraphs(1).Range.Find.Execute("Questions: The product rule", $false, $false, $false, $false, $false, $true, 1, $false, "Questions: The product rule", 2) | Out-Null
$d.Paragraphs(2).Range.Find.Execute("Sara Delgado Garcia", $false, $false, $false, $false, $false, $true, 1, $false, "Sara Delgado Garcia", 2) | Out-Null
$d.Paragraphs(4).Range.Find.Execute("A selection of questions for the study guide on the product rule.", $false, $false, $false, $false, $false, $true, 1, $false, "A selection of questions for the study guide on the product rule.", 2) | Out-Null

# --- Reorder m:sepChr before m:endChr inside every m:dPr (all OMath delimiters) ---
$om1 = $d.OMaths.Item(1)
$om1.Range.InsertXML('<m:oMathPara><m:oMathParaPr><m:jc m:val="center" /></m:oMathParaPr><m:oMath><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cosh</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:t>x</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>x</m:t></m:r></m:sup></m:sSup></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:nor /><m:sty m:val="p" /><m:scr m:val="sans-serif" /></m:rPr><m:t> and </m:t></m:r><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sinh</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:type m:val="bar" /></m:fPr><m:num><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:t>x</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>x</m:t></m:r></m:sup></m:sSup></m:num><m:den><m:r><m:t>2</m:t></m:r></m:den></m:f></m:oMath></m:oMathPara>')
$om4 = $d.OMaths.Item(4)
$om4.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:t>5</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>tan</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om5 = $d.OMaths.Item(5)
$om5.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om6 = $d.OMaths.Item(6)
$om6.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:d><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d></m:oMath>')
$om7 = $d.OMaths.Item(7)
$om7.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>13</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>5</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:d><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:d></m:oMath>')
$om8 = $d.OMaths.Item(8)
$om8.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:t>x</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>5</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:d><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>1</m:t></m:r></m:e></m:d></m:oMath>')
$om9 = $d.OMaths.Item(9)
$om9.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>10</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>21</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om10 = $d.OMaths.Item(10)
$om10.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cosh</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>2</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sinh</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om11 = $d.OMaths.Item(11)
$om11.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om12 = $d.OMaths.Item(12)
$om12.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:rad><m:radPr><m:degHide m:val="on" /></m:radPr><m:deg /><m:e><m:r><m:t>x</m:t></m:r></m:e></m:rad></m:oMath>')
$om13 = $d.OMaths.Item(13)
$om13.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cosh</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om14 = $d.OMaths.Item(14)
$om14.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:rad><m:radPr><m:degHide m:val="on" /></m:radPr><m:deg /><m:e><m:r><m:t>x</m:t></m:r></m:e></m:rad><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om15 = $d.OMaths.Item(15)
$om15.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>5</m:t></m:r><m:r><m:t>x</m:t></m:r></m:sup></m:sSup><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>3</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>5</m:t></m:r></m:e></m:d></m:oMath>')
$om16 = $d.OMaths.Item(16)
$om16.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>cos</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om17 = $d.OMaths.Item(17)
$om17.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>100</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om18 = $d.OMaths.Item(18)
$om18.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>5</m:t></m:r><m:r><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>2</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om19 = $d.OMaths.Item(19)
$om19.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>ln</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d></m:oMath>')
$om20 = $d.OMaths.Item(20)
$om20.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>5</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r></m:e></m:d><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>7</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:sSup><m:e><m:r><m:t>x</m:t></m:r></m:e><m:sup><m:r><m:t>4</m:t></m:r></m:sup></m:sSup></m:e></m:d></m:oMath>')
$om21 = $d.OMaths.Item(21)
$om21.Range.InsertXML('<m:oMath><m:r><m:t> </m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>sin</m:t></m:r><m:d><m:dPr><m:begChr m:val="(" /><m:sepChr m:val="" /><m:endChr m:val=")" /><m:grow /></m:dPr><m:e><m:r><m:t>x</m:t></m:r></m:e></m:d><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>+</m:t></m:r><m:r><m:t>3</m:t></m:r><m:r><m:t>x</m:t></m:r></m:e></m:d><m:sSup><m:e><m:r><m:t>e</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p" /></m:rPr><m:t>−</m:t></m:r><m:r><m:t>x</m:t></m:r></m:sup></m:sSup></m:oMath>')

Write-Output "done"
